$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 272
$ws.Range("C3").Value = 169511
$ws.Range("C4").Value = 160349
$ws.Range("C7").Value = 5.4
$ws.Range("C8").Value = 65.5
